# Generate Report for Handback
# Update the recorded Handoff/Handback timestamps for the
# 45a86069-05f5-4047-88b0-d30a99d4b967 file on both the zh-cn and de-de
# report sheets (row 2 of each sheet, columns E = "Correspond Handoff
# Datetime" and H = "Correspond Handback DateTime").

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-13 23:18:36"
$wsZh.Range("H2").Value = "2016-03-13 23:18:54"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-13 23:18:40"
$wsDe.Range("H2").Value = "2016-03-13 23:19:00"
